$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Export")

# Insert a new row above row 5 (which currently holds 004813166 / VENIA)
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "005681354"
$ws.Cells.Item(5, 1).ClearFormats()
$ws.Cells.Item(5, 2).Value = "MATHEUS"
$ws.Cells.Item(5, 3).Value = 65005.12
